# Germany Regionalliga North - update of league bases (19-06-2024 21:51)
#
# The underlying source data had a handful of match rows whose result /
# odds data had been attached to the wrong fixture. This script corrects
# that by exchanging the full data (every column except the running
# index in column A) between the affected rows:
#   - row 8  <-> row 9
#   - row 20 <-> row 21
#   - rows 95 -> 96 -> 97 -> 95 (three-way rotation)
#   - row 119 <-> row 120
#
# Columns B..AD hold: match id, competition, date, home team, away team,
# full/half time score, result, odds, etc. Column A (running number) is
# left untouched because it simply reflects the row position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data columns that can be exchanged between two fixture rows
# (everything except column A).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Get-RowValues($row) {
    $values = @{}
    foreach ($col in $cols) {
        $values[$col] = $ws.Range("$col$row").Value2
    }
    return $values
}

function Set-RowValues($row, $values) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

# --- Simple swaps (two fixtures exchanged with each other) ---
function Swap-Rows($rowA, $rowB) {
    $valsA = Get-RowValues $rowA
    $valsB = Get-RowValues $rowB
    Set-RowValues $rowA $valsB
    Set-RowValues $rowB $valsA
}

Swap-Rows 8 9
Swap-Rows 20 21
Swap-Rows 119 120

# --- Three-way rotation: new(95) = old(96), new(96) = old(97), new(97) = old(95) ---
$vals95 = Get-RowValues 95
$vals96 = Get-RowValues 96
$vals97 = Get-RowValues 97

Set-RowValues 95 $vals96
Set-RowValues 96 $vals97
Set-RowValues 97 $vals95
